# Update cryptocurrency price/volume table with latest values.
# Generated from the target diff: updates Price (D) and Volume(1h) (E) columns
# for all rows, plus Coin name (B) and Link (C) for rows whose ranking order
# changed (18/19, 31/32, 33/34, 43/44, 48/49).
#
# D-column price values are forced to text (leading apostrophe) so Excel's
# automatic type detection does not convert them to numbers, which would
# silently strip meaningful trailing zeros (e.g. "1.00" -> 1) or reformat
# values using a different decimal representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.791.76"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").Value = "3.481.83"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'409.73"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'132.09"
$ws.Range("E6").Value = "  +17.28%  "
$ws.Range("D7").Value = "3.475.64"
$ws.Range("E7").Value = "  +4.01%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.696"
$ws.Range("E10").Value = "  +9.46%  "
$ws.Range("D11").Value = "'0.130"
$ws.Range("E11").Value = "  +30.90%  "
$ws.Range("D12").Value = "'43.22"
$ws.Range("E12").Value = "  +7.48%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "4.007.66"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "'8.78"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "'20.32"
$ws.Range("E16").Value = "  +4.77%  "
$ws.Range("D17").Value = "3.469.29"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.06"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "62.644.11"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").Value = "'10.96"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "'0.0000139"
$ws.Range("E21").Value = "  +24.93%  "
$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "'82.94"
$ws.Range("E23").Value = "  +9.23%  "
$ws.Range("D24").Value = "'13.22"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'312.71"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "'3.19"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'30.53"
$ws.Range("E27").Value = "  +6.39%  "
$ws.Range("D28").Value = "'8.35"
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("D29").Value = "'7.79"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("D30").Value = "'0.180"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +4.75%  "
$ws.Range("B32").Value = "LEO"
$ws.Range("C32").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D32").Value = "'4.38"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.98"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.67"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("D35").Value = "'43.38"
$ws.Range("E35").Value = "  +8.53%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'0.0496"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "'52.62"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +6.70%  "
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").Value = "'0.127"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'2.00"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'138.25"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "'17.70"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("D46").Value = "'4.02"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'0.288"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'22.69"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.25"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "2.210.63"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "3.814.00"
$ws.Range("E51").Value = "  +3.88%  "

